$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 95,3
$data[0,0] = "Cluster name"
$data[0,1] = "Active cases"
$data[0,2] = "Exist"
$data[1,0] = "253 Hoddle Street Apartment Complex Collingwood"
$data[1,1] = 5
$data[1,2] = "old"
$data[2,0] = "3535 Opal Meadow Heights Aged Care Community Meadow Heights"
$data[2,1] = 27
$data[2,2] = "old"
$data[3,0] = "3535 Opal Meadow Heights Aged Care Community Meadow Heights"
$data[3,1] = 28
$data[3,2] = "new"
$data[4,0] = "Adorn Cosmetics Clayton"
$data[4,1] = 5
$data[4,2] = "new"
$data[5,0] = "Al Haj Halal Meats Glenroy"
$data[5,1] = 34
$data[5,2] = "new"
$data[6,0] = "Al Haj Halal Meats Glenroy"
$data[6,1] = 37
$data[6,2] = "old"
$data[7,0] = "Al-Taqwa College Truganina"
$data[7,1] = 5
$data[7,2] = "new"
$data[8,0] = "Al-Taqwa College Truganina"
$data[8,1] = 6
$data[8,2] = "old"
$data[9,0] = "Amiga Montessori Craigieburn"
$data[9,1] = 25
$data[9,2] = "old"
$data[10,0] = "CS Square Caroline Springs"
$data[10,1] = 7
$data[10,2] = "old"
$data[11,0] = "CS Square Caroline Springs"
$data[11,1] = 9
$data[11,2] = "new"
$data[12,0] = "Cannie Road Construction Site Cannie"
$data[12,1] = 7
$data[12,2] = "new"
$data[13,0] = "Cannie Road Construction Site Cannie"
$data[13,1] = 8
$data[13,2] = "old"
$data[14,0] = "Cedars Medical Clinic Coburg"
$data[14,1] = 28
$data[14,2] = "new"
$data[15,0] = "Cedars Medical Clinic Coburg"
$data[15,1] = 41
$data[15,2] = "new"
$data[16,0] = "Cedars Medical Clinic Coburg"
$data[16,1] = 42
$data[16,2] = "old"
$data[17,0] = "Chemist Warehouse Campbellfield DC"
$data[17,1] = 5
$data[17,2] = "new"
$data[18,0] = "City of Hobsons Bay Community"
$data[18,1] = 5
$data[18,2] = "old"
$data[19,0] = "City of Moreland Community"
$data[19,1] = 6
$data[19,2] = "new"
$data[20,0] = "City of Moreland Community"
$data[20,1] = 7
$data[20,2] = "old"
$data[21,0] = "City of Wyndham Community"
$data[21,1] = 6
$data[21,2] = "new"
$data[22,0] = "City of Wyndham Community"
$data[22,1] = 7
$data[22,2] = "old"
$data[23,0] = "Coles Campbellfield Plaza Campbellfield"
$data[23,1] = 8
$data[23,2] = "old"
$data[24,0] = "Coles Campbellfield Plaza Campbellfield"
$data[24,1] = 9
$data[24,2] = "new"
$data[25,0] = "Coles Coburg North Village"
$data[25,1] = 26
$data[25,2] = "new"
$data[26,0] = "Coles Coburg North Village"
$data[26,1] = 27
$data[26,2] = "old"
$data[27,0] = "Coles Pakenham Place Shopping Centre"
$data[27,1] = 10
$data[27,2] = "old"
$data[28,0] = "Coles Pakenham Place Shopping Centre"
$data[28,1] = 11
$data[28,2] = "new"
$data[29,0] = "Coles Roxburgh Village Roxburgh Park"
$data[29,1] = 9
$data[29,2] = "old"
$data[30,0] = "Coles Roxburgh Village Roxburgh Park"
$data[30,1] = 11
$data[30,2] = "new"
$data[31,0] = "Community Kids Meadow Heights"
$data[31,1] = 12
$data[31,2] = "new"
$data[32,0] = "Community Kids Meadow Heights"
$data[32,1] = 14
$data[32,2] = "old"
$data[33,0] = "Construction Site Olea Apartment Caulfield North"
$data[33,1] = 12
$data[33,2] = "old"
$data[34,0] = "Construction Site Olea Apartment Caulfield North"
$data[34,1] = 16
$data[34,2] = "new"
$data[35,0] = "Costco Wholesale Epping"
$data[35,1] = 18
$data[35,2] = "old"
$data[36,0] = "Costco Wholesale Epping"
$data[36,1] = 24
$data[36,2] = "new"
$data[37,0] = "Crusader Caravans Epping"
$data[37,1] = 13
$data[37,2] = "old"
$data[38,0] = "Crusader Caravans Epping"
$data[38,1] = 14
$data[38,2] = "new"
$data[39,0] = "Direct Freight Express Cambellfield"
$data[39,1] = 13
$data[39,2] = "new"
$data[40,0] = "Direct Freight Express Campbellfield"
$data[40,1] = 13
$data[40,2] = "old"
$data[41,0] = "Epworth Healthcare Epworth Richmond Emergency Department"
$data[41,1] = 5
$data[41,2] = "new"
$data[42,0] = "Fitzroy Community School Fitzroy North"
$data[42,1] = 36
$data[42,2] = "new"
$data[43,0] = "Fitzroy Community School Fitzroy North"
$data[43,1] = 41
$data[43,2] = "old"
$data[44,0] = "Glenroy West Primary School"
$data[44,1] = 5
$data[44,2] = "old"
$data[45,0] = "Glenroy West Primary School"
$data[45,1] = 6
$data[45,2] = "new"
$data[46,0] = "Goodstart Early Learning Altona"
$data[46,1] = 5
$data[46,2] = "new"
$data[47,0] = "Green Leaves Early Learning Centre Highlands Craigieburn"
$data[47,1] = 7
$data[47,2] = "new"
$data[48,0] = "IGA Meadow Heights Shopping Centre Meadow Heights"
$data[48,1] = 6
$data[48,2] = "new"
$data[49,0] = "IGA Meadow Heights Shopping Centre Meadow Heights"
$data[49,1] = 7
$data[49,2] = "old"
$data[50,0] = "ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine"
$data[50,1] = 7
$data[50,2] = "old"
$data[51,0] = "ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine"
$data[51,1] = 9
$data[51,2] = "new"
$data[52,0] = "Ilim College Glenroy Campus Hadfield"
$data[52,1] = 16
$data[52,2] = "new"
$data[53,0] = "Ilim College Glenroy Campus Hadfield"
$data[53,1] = 19
$data[53,2] = "old"
$data[54,0] = "Ilim Learning Sanctuary Glenroy"
$data[54,1] = 11
$data[54,2] = "new"
$data[55,0] = "Ilim Learning Sanctuary Glenroy"
$data[55,1] = 12
$data[55,2] = "old"
$data[56,0] = "Impact Designer Homes Epping"
$data[56,1] = 5
$data[56,2] = "new"
$data[57,0] = "Industrial Galvanizers Valmont Coatings Campbellfield"
$data[57,1] = 18
$data[57,2] = "old"
$data[58,0] = "Industrial Galvanizers Valmont Coatings Campbellfield"
$data[58,1] = 22
$data[58,2] = "new"
$data[59,0] = "Islamic College of Melbourne Tarneit"
$data[59,1] = 5
$data[59,2] = "new"
$data[60,0] = "Islamic College of Melbourne Tarneit"
$data[60,1] = 9
$data[60,2] = "old"
$data[61,0] = "Kasr Sweets Coolaroo"
$data[61,1] = 5
$data[61,2] = "new"
$data[62,0] = "Kasr Sweets Coolaroo"
$data[62,1] = 6
$data[62,2] = "old"
$data[63,0] = "Learning Nest Early Learning Centre Meadow Heights"
$data[63,1] = 5
$data[63,2] = "new"
$data[64,0] = "Learning Nest Early Learning Centre Meadow Heights"
$data[64,1] = 6
$data[64,2] = "old"
$data[65,0] = "MyCentre Childcare Broadmeadows"
$data[65,1] = 14
$data[65,2] = "old"
$data[66,0] = "MyCentre Childcare Broadmeadows"
$data[66,1] = 17
$data[66,2] = "new"
$data[67,0] = "Newbury Child and Community Centre Craigieburn"
$data[67,1] = 5
$data[67,2] = "new"
$data[68,0] = "Newbury Child and Community Centre Craigieburn"
$data[68,1] = 7
$data[68,2] = "old"
$data[69,0] = "Nino Early Learning Adventures Lalor"
$data[69,1] = 5
$data[69,2] = "old"
$data[70,0] = "Northern Health Northern Hospital Epping Emergency Department Tier 1B"
$data[70,1] = 44
$data[70,2] = "old"
$data[71,0] = "Northern Health Northern Hospital Epping Emergency Department Tier 1B"
$data[71,1] = 47
$data[71,2] = "new"
$data[72,0] = "Northern Health The Northern Hospital Epping"
$data[72,1] = 5
$data[72,2] = "new"
$data[73,0] = "Our Lady Help of Christian's Primary School Brunswick East"
$data[73,1] = 8
$data[73,2] = "new"
$data[74,0] = "Paisley Park Early Learning Centre Bundoora"
$data[74,1] = 8
$data[74,2] = "new"
$data[75,0] = "Paisley Park Early Learning Centre Bundoora"
$data[75,1] = 9
$data[75,2] = "old"
$data[76,0] = "Panorama Construction Site Whitehorse Rd Box Hill"
$data[76,1] = 17
$data[76,2] = "new"
$data[77,0] = "Panorama Construction Site Whitehorse Rd Box Hill"
$data[77,1] = 20
$data[77,2] = "old"
$data[78,0] = "People First Healthcare Home Residence Disability Support Taylors Lakes"
$data[78,1] = 5
$data[78,2] = "old"
$data[79,0] = "Private Residence Northern Community Services Fawkner"
$data[79,1] = 5
$data[79,2] = "new"
$data[80,0] = "Richmond Quarter 261-271 Bridge Road Construction Site Richmond"
$data[80,1] = 9
$data[80,2] = "old"
$data[81,0] = "Richmond Quarter 261-271 Bridge Road Construction Site Richmond"
$data[81,1] = 10
$data[81,2] = "new"
$data[82,0] = "Sharpline Stainless Steel Coburg North"
$data[82,1] = 5
$data[82,2] = "new"
$data[83,0] = "Tek Foods Somerton"
$data[83,1] = 12
$data[83,2] = "new"
$data[84,0] = "Tek Foods Somerton"
$data[84,1] = 13
$data[84,2] = "old"
$data[85,0] = "The Homestead Child and Family Centre Roxburgh Park"
$data[85,1] = 12
$data[85,2] = "new"
$data[86,0] = "The Homestead Child and Family Centre Roxburgh Park"
$data[86,1] = 13
$data[86,2] = "old"
$data[87,0] = "The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B"
$data[87,1] = 10
$data[87,2] = "old"
$data[88,0] = "The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B"
$data[88,1] = 11
$data[88,2] = "new"
$data[89,0] = "Western Health Footscray Hospital Emergency Department"
$data[89,1] = 5
$data[89,2] = "old"
$data[90,0] = "Western Health Footscray Hospital Emergency Department"
$data[90,1] = 6
$data[90,2] = "new"
$data[91,0] = "Western Health Sunshine Hospital Emergency Department"
$data[91,1] = 7
$data[91,2] = "old"
$data[92,0] = "Western Health Sunshine Hospital Emergency Department"
$data[92,1] = 8
$data[92,2] = "new"
$data[93,0] = "Woolworths Greenvale Lakes Roxburgh Park"
$data[93,1] = 5
$data[93,2] = "new"
$data[94,0] = "Woolworths Greenvale Lakes Roxburgh Park"
$data[94,1] = 6
$data[94,2] = "old"

$ws.Range("A1:C95").Value = $data

# Clear any leftover rows from the previous (longer) table
$ws.Range("A96:C103").ClearContents()